$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column header in H1, copying the same formatting
# used by the other header cells (e.g. G1 - bold, bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the "Save" values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
